$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Clear rows/cells that are removed entirely in the target layout ---
$ws.Range("G47:H47").ClearContents()
$ws.Range("G49:H49").ClearContents()
$ws.Range("G51:H51").ClearContents()
$ws.Range("G22:H22").ClearContents()
$ws.Range("H25").ClearContents()
$ws.Range("H30").ClearContents()
$ws.Range("G38").ClearContents()

# --- Prepare new rows 37 and 38 with same formatting as row 36 ---
$ws.Range("B36:C36").Copy()
$ws.Range("B37:C37").PasteSpecial(-4122)
$ws.Range("B36:C36").Copy()
$ws.Range("B38:C38").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Set cell values to match target state ---
$ws.Range('B2').Value = 'Date'
$ws.Range('C2').Value = 'Action'
$ws.Range('G2').Value = 'Feature'
$ws.Range('H2').Value = 'Backlog'
$ws.Range('I2').Value = 'Priority'
$ws.Range('B4').Value = '2/7/2017'
$ws.Range('C4').Value = 'Refactored carService to return promise (via SO question)'
$ws.Range('G4').Value = 'Wishlist'
$ws.Range('H4').Value = 'Add Wishlist Table - AspNetUser Id and Car ID'
$ws.Range('J4').Value = 'Search DDL''s - only use values when $dirty flagged'
$ws.Range('B5').Value = '2/7/2017'
$ws.Range('C5').Value = 'Render pagination links via ng-repeat'
$ws.Range('G5').Value = 'WishList'
$ws.Range('H5').Value = 'Ad to Wishlist - adds record to Db'
$ws.Range('J5').Value = 'Populate DDLS from API calls'
$ws.Range('B6').Value = '2/7/2017'
$ws.Range('C6').Value = 'Set active class of current page when clicked on pagination links.'
$ws.Range('G6').Value = 'Wishlist'
$ws.Range('H6').Value = 'Add Heart & Heart-o - When users adds it chage fa icon'
$ws.Range('J6').Value = 'Filter DDLs by other selections?'
$ws.Range('B7').Value = '2/8/2017'
$ws.Range('C7').Value = 'Moved Categories partial view to angular directive'
$ws.Range('G7').Value = 'Wishlist'
$ws.Range('H7').Value = 'Hook up link to show Wishlist and display'
$ws.Range('B8').Value = '2/8/2017'
$ws.Range('C8').Value = 'Investigate Pagination bug & fix'
$ws.Range('G8').Value = 'Wishlist'
$ws.Range('H8').Value = 'Remove from WishList - remove from Db'
$ws.Range('B9').Value = '2/8/2017'
$ws.Range('C9').Value = 'Make start on Details page - pass in Id value'
$ws.Range('B10').Value = '2/8/2017'
$ws.Range('C10').Value = 'Answer 2 x Stack Overflow Questions'
$ws.Range('G10').Value = 'Login'
$ws.Range('H10').Value = 'User can log in from modal'
$ws.Range('B11').Value = '2/8/2017'
$ws.Range('C11').Value = 'Read chapter of Web API 2 book'
$ws.Range('B12').Value = '2/9/2017'
$ws.Range('C12').Value = 'Get all Car Info Data (Console app - wikiApi)'
$ws.Range('G12').Value = 'Homepage'
$ws.Range('H12').Value = 'Finish update - complete text, fa icons etc'
$ws.Range('B13').Value = '2/9/2017'
$ws.Range('C13').Value = 'Added car spec table directive to car details page'
$ws.Range('G13').Value = 'Homepage'
$ws.Range('H13').Value = 'Add Google chart to Homepage'
$ws.Range('B14').Value = '2/9/2017'
$ws.Range('C14').Value = 'Read chapter of Design Patterns Book'
$ws.Range('B15').Value = '2/10/2017'
$ws.Range('C15').Value = 'Added Car Info to DB - added car info service and output to details page'
$ws.Range('G15').Value = 'Recommended'
$ws.Range('H15').Value = 'Recommended & similar cars (dynamic from api controller) Details Page'
$ws.Range('B16').Value = '2/10/2017'
$ws.Range('C16').Value = 'Add images from API call to carousel on Details page'
$ws.Range('B17').Value = '2/11/2017'
$ws.Range('C17').Value = 'Added search algorithm'
$ws.Range('G17').Value = 'Reviews'
$ws.Range('H17').Value = 'Reviews - and Ratings - Stars fa icons'
$ws.Range('B18').Value = '2/11/2017'
$ws.Range('C18').Value = 'Investigate Partial view, directive issue - Header'
$ws.Range('B19').Value = '2/12/2017'
$ws.Range('C19').Value = 'Hooked up text search to front end - full text search'
$ws.Range('G19').Value = 'Checkout'
$ws.Range('H19').Value = 'Make checkout & confirmation pages dynamic'
$ws.Range('B20').Value = '2/12/2017'
$ws.Range('C20').Value = 'Categories CSS classes toggle'
$ws.Range('G20').Value = 'Checkout'
$ws.Range('H20').Value = 'Email With Confirmation - and PDF Invoice'
$ws.Range('B21').Value = '2/12/2017'
$ws.Range('C21').Value = 'Created Backlog with all outstanding tasks'
$ws.Range('G21').Value = 'Checkout'
$ws.Range('H21').Value = 'Cost Of Rental and calculate based on special offers etc - Update Cost to rent per week'
$ws.Range('B22').Value = '2/13/2017'
$ws.Range('C22').Value = 'Unit Tests for ar Repo and TextSearch'
$ws.Range('B23').Value = '2/13/2017'
$ws.Range('C23').Value = 'Updated Home page with pics & text & video'
$ws.Range('G23').Value = 'Special Offer'
$ws.Range('H23').Value = 'Add Bool for Featured Cars, Cars On special '
$ws.Range('B24').Value = '2/14/2017'
$ws.Range('C24').Value = 'Fixed Header partial view issue - now calls search'
$ws.Range('G24').Value = 'Special Offer'
$ws.Range('H24').Value = 'link to special offer/ featuredlist on click Call To Action Links'
$ws.Range('B25').Value = '2/14/2017'
$ws.Range('C25').Value = 'fragaria/angular-daterangepicker - not working due to jquery issue'
$ws.Range('B26').Value = '2/15/2017'
$ws.Range('C26').Value = 'Added more unit tests for text search'
$ws.Range('H26').Value = 'Daterange picker - values, dates to angular controller'
$ws.Range('B27').Value = '2/16/2017'
$ws.Range('C27').Value = 'Fade in - scroll directive opacity - added'
$ws.Range('H27').Value = 'Finish Categories Sidebar'
$ws.Range('B28').Value = '2/16/2017'
$ws.Range('C28').Value = 'Added and updated search Pagination'
$ws.Range('H28').Value = 'Owl Carousel data-img-zoom not working'
$ws.Range('B29').Value = '2/16/2017'
$ws.Range('C29').Value = 'Added Recommended Rentals to Details page - static'
$ws.Range('H29').Value = 'Footer links'
$ws.Range('B30').Value = '2/16/2017'
$ws.Range('C30').Value = 'Add appstrap shop-checkout & confirmation pages static'
$ws.Range('B31').Value = '2/16/2017'
$ws.Range('C31').Value = 'Added new date picker - displays start and end dates'
$ws.Range('H31').Value = 'Return IHttpActionResult from apiControllers  - Ok, Bad Result etc.'
$ws.Range('B32').Value = '2/17/2017'
$ws.Range('C32').Value = 'Pagination page links - SO question - (created Plunker)'
$ws.Range('H32').Value = 'Unit Tests - Nunit'
$ws.Range('B33').Value = '2/17/2017'
$ws.Range('C33').Value = 'Added wishlistController and can get carId''s from link click - Details & List'
$ws.Range('H33').Value = 'Jasmine Tests JS '
$ws.Range('B34').Value = '2/18/2017'
$ws.Range('C34').Value = 'Added Log In and Register pages - static'
$ws.Range('H34').Value = 'AutoComplete on text boxes???'
$ws.Range('B35').Value = '2/19/2017'
$ws.Range('C35').Value = 'Register working - user added to AspNetUser table'
$ws.Range('B36').Value = '2/21/2017'
$ws.Range('C36').Value = 'Checkout populated from View Model'
$ws.Range('H36').Value = 'Check Thumbnail pics for cars - and change ones not adequate'
$ws.Range('B37').Value = '2/21/2017'
$ws.Range('C37').Value = 'Checkout & Confirmation Pages Working'
$ws.Range('B38').Value = '2/22/2017'
$ws.Range('C38').Value = 'Change page size, number results, e.g. 15, 20 etc.'
$ws.Range('H38').Value = 'Add Car List as Grid View'
$ws.Range('G40').Value = 'Search'
$ws.Range('H40').Value = 'Add Search Bar - 3 x Drop Down Lists - filtered by each selection, make, model, year'
$ws.Range('H42').Value = 'Trending Now'
$ws.Range('G44').Value = 'Validation'
$ws.Range('H44').Value = 'Validation logic and add modelstate errors etc.'
$ws.Range('G46').Value = 'Pagination'
$ws.Range('H46').Value = 'Go direct to certain page - eg pg 17'
$ws.Range('G48').Value = 'Logging'
$ws.Range('H48').Value = 'Log errors with Nlog?'
$ws.Range('G50').Value = 'Admin Area'
$ws.Range('H50').Value = 'Update cars etc.'

# --- Misc formatting updates ---
$ws.Columns(9).ColumnWidth = 6.736979166666667
$ws.Columns(10).ColumnWidth = 44.592447916666664
$ws.Range("H53").Select()
